$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$valsG = @(
    "assertAttribute(locator,attribute,expected)",
    "assertChecked(name)",
    "assertDisabled(name)",
    "assertElementPresent(name)",
    "assertEnabled(name)",
    "assertHierCells(matchBy,column,expected,nestedOnly)",
    "assertHierRow(matchBy,expected)",
    "assertListCount(count)",
    "assertLocatorNotPresent(locator)",
    "assertLocatorPresent(locator)",
    "assertMenuEnabled(menu)",
    "assertModalDialogNotPresent()",
    "assertModalDialogPresent()",
    "assertModalDialogTitle(title)",
    "assertModalDialogTitleByLocator(locator,title)",
    "assertNotChecked(name)",
    "assertSelected(name,text)",
    "assertTableCell(row,column,contains)",
    "assertTableColumnContains(column,contains)",
    "assertTableContains(contains)",
    "assertTableRowContains(row,contains)",
    "assertText(name,expected)",
    "assertWindowTitleContains(contains)",
    "clear(locator)",
    "clearCombo(name)",
    "clearModalDialog(var,button)",
    "clearTextArea(name)",
    "clearTextBox(name)",
    "clickButton(name)",
    "clickByLocator(locator)",
    "clickCheckBox(name)",
    "clickElementOffset(name,xOffset,yOffset)",
    "clickExplorerBar(group,item)",
    "clickFirstMatchRow(nameValues)",
    "clickFirstMatchedList(contains)",
    "clickIcon(label)",
    "clickList(row)",
    "clickMenu(menu)",
    "clickOffset(locator,xOffset,yOffset)",
    "clickRadio(name)",
    "clickTab(group,name)",
    "clickTableCell(row,column)",
    "clickTableRow(row)",
    "clickTextPane(name,criteria)",
    "clickTextPaneRow(var,index)",
    "closeApplication()",
    "collapseHierTable()",
    "editCurrentRow(nameValues)",
    "editHierCells(var,matchBy,nameValues)",
    "editTableCells(row,nameValues)",
    "getRowCount(var)",
    "hideExplorerBar()",
    "login(form,username,password)",
    "maximize()",
    "minimize()",
    "resize(width,height)",
    "saveAllTableRows(var)",
    "saveAttributeByLocator(var,locator,attribute)",
    "saveElementCount(var,name)",
    "saveFirstListData(var,contains)",
    "saveFirstMatchedListIndex(var,contains)",
    "saveHierCells(var,matchBy,column,nestedOnly)",
    "saveHierRow(var,matchBy)",
    "saveListData(var,contains)",
    "saveLocatorCount(var,locator)",
    "saveModalDialogText(var)",
    "saveModalDialogTextByLocator(var,locator)",
    "saveProcessId(var,locator)",
    "saveRowCount(var)",
    "saveTableRows(var,contains)",
    "saveTableRowsRange(var,beginRow,endRow)",
    "saveText(var,name)",
    "saveTextByLocator(var,locator)",
    "saveTextPane(var,name,criteria)",
    "saveWindowTitle(var)",
    "scanTable(var,name)",
    "selectCombo(name,text)",
    "sendKeysToTextBox(name,text1,text2,text3,text4)",
    "showExplorerBar()",
    "toggleExplorerBar()",
    "typeAppendTextArea(name,text1,text2,text3,text4)",
    "typeAppendTextBox(name,text1,text2,text3,text4)",
    "typeByLocator(locator,text)",
    "typeKeys(os,keystrokes)",
    "typeTextArea(name,text1,text2,text3,text4)",
    "typeTextBox(name,text1,text2,text3,text4)",
    "useApp(appId)",
    "useForm(formName)",
    "useHierTable(var,name)",
    "useList(var,name)",
    "useTable(var,name)",
    "useTableRow(var,row)",
    "waitFor(name,maxWaitMs)",
    "waitForLocator(locator,maxWaitMs)"
)
for ($i = 0; $i -lt $valsG.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $valsG[$i]
}

$valsJ = @(
    "colorbit(source,bit,saveTo)",
    "compare(baseline,actual)",
    "convert(source,format,saveTo)",
    "crop(image,dimension,saveTo)",
    "resize(image,width,height,saveTo)",
    "saveDiff(var,baseline,actual)"
)
for ($i = 0; $i -lt $valsJ.Length; $i++) {
    $ws.Cells.Item($i + 2, 10).Value = $valsJ[$i]
}

$valsM = @(
    "addOrReplace(json,jsonpath,input,var)",
    "assertCorrectness(json,schema)",
    "assertElementCount(json,jsonpath,count)",
    "assertElementNotPresent(json,jsonpath)",
    "assertElementPresent(json,jsonpath)",
    "assertEqual(expected,actual)",
    "assertValue(json,jsonpath,expected)",
    "assertValues(json,jsonpath,array,exactOrder)",
    "assertWellformed(json)",
    "beautify(json,var)",
    "compact(var,json,removeEmpty)",
    "fromCsv(csv,header,jsonFile)",
    "minify(json,var)",
    "storeCount(json,jsonpath,var)",
    "storeValue(json,jsonpath,var)",
    "storeValues(json,jsonpath,var)"
)
for ($i = 0; $i -lt $valsM.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $valsM[$i]
}

$valsY = @(
    "assertAndClick(locator,label)",
    "assertAttribute(locator,attrName,value)",
    "assertAttributeContains(locator,attrName,contains)",
    "assertAttributeNotContains(locator,attrName,contains)",
    "assertAttributeNotPresent(locator,attrName)",
    "assertAttributePresent(locator,attrName)",
    "assertChecked(locator)",
    "assertContainCount(locator,text,count)",
    "assertCssNotPresent(locator,property)",
    "assertCssPresent(locator,property,value)",
    "assertElementByAttributes(nameValues)",
    "assertElementByText(locator,text)",
    "assertElementCount(locator,count)",
    "assertElementNotPresent(locator)",
    "assertElementPresent(locator)",
    "assertElementsPresent(prefix)",
    "assertFocus(locator)",
    "assertFrameCount(count)",
    "assertFramePresent(frameName)",
    "assertIECompatMode()",
    "assertIENativeMode()",
    "assertLinkByLabel(label)",
    "assertNotChecked(locator)",
    "assertNotFocus(locator)",
    "assertNotText(locator,text)",
    "assertNotVisible(locator)",
    "assertOneMatch(locator)",
    "assertScrollbarHNotPresent(locator)",
    "assertScrollbarHPresent(locator)",
    "assertScrollbarVNotPresent(locator)",
    "assertScrollbarVPresent(locator)",
    "assertTable(locator,row,column,text)",
    "assertText(locator,text)",
    "assertTextContains(locator,text)",
    "assertTextCount(locator,text,count)",
    "assertTextList(locator,list,ignoreOrder)",
    "assertTextMatches(text,minMatch,scrollTo)",
    "assertTextNotContains(locator,text)",
    "assertTextNotPresent(text)",
    "assertTextOrder(locator,descending)",
    "assertTextPresent(text)",
    "assertTitle(text)",
    "assertValue(locator,value)",
    "assertValueOrder(locator,descending)",
    "assertVisible(locator)",
    "checkAll(locator)",
    "clearLocalStorage()",
    "click(locator)",
    "clickAll(locator)",
    "clickAndWait(locator,waitMs)",
    "clickByLabel(label)",
    "clickByLabelAndWait(label,waitMs)",
    "clickOffset(locator,x,y)",
    "clickWithKeys(locator,keys)",
    "close()",
    "closeAll()",
    "deselect(locator,text)",
    "deselectMulti(locator,array)",
    "dismissInvalidCert()",
    "dismissInvalidCertPopup()",
    "doubleClick(locator)",
    "doubleClickAndWait(locator,waitMs)",
    "doubleClickByLabel(label)",
    "doubleClickByLabelAndWait(label,waitMs)",
    "dragAndDrop(fromLocator,toLocator)",
    "dragTo(fromLocator,xOffset,yOffset)",
    "editLocalStorage(key,value)",
    "executeScript(var,script)",
    "focus(locator)",
    "goBack()",
    "goBackAndWait()",
    "maximizeWindow()",
    "mouseOver(locator)",
    "open(url)",
    "openAndWait(url,waitMs)",
    "openHttpBasic(url,username,password)",
    "openIgnoreTimeout(url)",
    "refresh()",
    "refreshAndWait()",
    "resizeWindow(width,height)",
    "rightClick(locator)",
    "saveAllWindowIds(var)",
    "saveAllWindowNames(var)",
    "saveAttribute(var,locator,attrName)",
    "saveAttributeList(var,locator,attrName)",
    "saveCount(var,locator)",
    "saveDivsAsCsv(headers,rows,cells,nextPage,file)",
    "saveElement(var,locator)",
    "saveElements(var,locator)",
    "saveLocalStorage(var,key)",
    "saveLocation(var)",
    "savePageAs(var,sessionIdName,url)",
    "savePageAsFile(sessionIdName,url,file)",
    "saveTableAsCsv(locator,nextPageLocator,file)",
    "saveText(var,locator)",
    "saveTextArray(var,locator)",
    "saveTextSubstringAfter(var,locator,delim)",
    "saveTextSubstringBefore(var,locator,delim)",
    "saveTextSubstringBetween(var,locator,start,end)",
    "saveValue(var,locator)",
    "saveValues(var,locator)",
    "scrollElement(locator,xOffset,yOffset)",
    "scrollLeft(locator,pixel)",
    "scrollPage(xOffset,yOffset)",
    "scrollRight(locator,pixel)",
    "scrollTo(locator)",
    "select(locator,text)",
    "selectFrame(locator)",
    "selectMulti(locator,array)",
    "selectMultiOptions(locator)",
    "selectText(locator)",
    "selectWindow(winId)",
    "selectWindowAndWait(winId,waitMs)",
    "selectWindowByIndex(index)",
    "selectWindowByIndexAndWait(index,waitMs)",
    "toggleSelections(locator)",
    "type(locator,value)",
    "typeKeys(locator,value)",
    "uncheckAll(locator)",
    "unselectAllText()",
    "upload(fieldLocator,file)",
    "verifyContainText(locator,text)",
    "verifyText(locator,text)",
    "wait(waitMs)",
    "waitForElementPresent(locator)",
    "waitForPopUp(winId,waitMs)",
    "waitForTextPresent(text)",
    "waitForTitle(text)"
)
for ($i = 0; $i -lt $valsY.Length; $i++) {
    $ws.Cells.Item($i + 2, 25).Value = $valsY[$i]
}

$valsAD = @(
    "append(xml,xpath,content,var)",
    "assertCorrectness(xml,schema)",
    "assertElementCount(xml,xpath,count)",
    "assertElementNotPresent(xml,xpath)",
    "assertElementPresent(xml,xpath)",
    "assertSoap(wsdl,xml)",
    "assertSoapFaultCode(expected,xml)",
    "assertSoapFaultString(expected,xml)",
    "assertValue(xml,xpath,expected)",
    "assertValues(xml,xpath,array,exactOrder)",
    "assertWellformed(xml)",
    "beautify(xml,var)",
    "clear(xml,xpath,var)",
    "delete(xml,xpath,var)",
    "insertAfter(xml,xpath,content,var)",
    "insertBefore(xml,xpath,content,var)",
    "minify(xml,var)",
    "prepend(xml,xpath,content,var)",
    "replace(xml,xpath,content,var)",
    "replaceIn(xml,xpath,content,var)",
    "storeCount(xml,xpath,var)",
    "storeSoapFaultCode(var,xml)",
    "storeSoapFaultDetail(var,xml)",
    "storeSoapFaultString(var,xml)",
    "storeValue(xml,xpath,var)",
    "storeValues(xml,xpath,var)"
)
for ($i = 0; $i -lt $valsAD.Length; $i++) {
    $ws.Cells.Item($i + 2, 30).Value = $valsAD[$i]
}

$wb.Names.Item("desktop").RefersTo = "='#system'!`$G`$2:`$G`$95"
$wb.Names.Item("image").RefersTo = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$17"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"